$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: add a bottom border under the (now not-last) row -----------
# A4/B4 keep the normal (8-pt / size-0) font but gain a thin bottom border.
$ws.Range("A4:B4").Borders.Item(9).LineStyle = 1

# C4:D4:E4 keep their small font (size 8) and also gain the thin bottom
# border.
$rowRight = $ws.Range("C4:E4")
$rowRight.Font.Size = 8
$rowRight.Borders.Item(9).LineStyle = 1

# --- Row 5: new data row -------------------------------------------------
$ws.Range("B5").Value2 = 253
$ws.Range("C5").Value2 = ' It\''d be good to see that [CS:N]Grovyle[CR]\ncaught soon!'
$ws.Range("D5").Value2 = ' Приятно знать, что [CS:N]Гровайл[CR]\nскоро будет пойман!'
$ws.Range("E5").Value2 = ' Ðñéÿóîï èîàóû, œóï [CS:N]Ãñïâàêì[CR]\nòëïñï áôäåó ðïêíàî!'

# Match row 4's row height for the new row (auto-fit wrapped two-line text).
$ws.Rows.Item(5).RowHeight = 21.6

# --- Selection -------------------------------------------------------------
$ws.Range("D3").Select() | Out-Null
